# Apply "Score" column (D) data entry to Sheet1 of the marine mammals
# ESA status workbook, plus the trailing AVERAGE summary row and minor
# view-state tweaks, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header D1 from "IUCN status" to "Score".
$ws.Range("D1").Value = "Score"

# Row -> Score value, as collected in the source data.
$scores = @{
    2  = 0
    3  = 0.6
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0.6
    9  = 0.6
    10 = 0
    11 = 0.5
    12 = 0.6
    13 = 0.6
    14 = 0.6
    15 = 0
    16 = 0.6
    17 = 0.6
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0.5
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0.6
    29 = 0.6
    30 = 0.6
    31 = 0
    32 = 0
}

foreach ($r in $scores.Keys) {
    $ws.Cells.Item($r, 4).Value = $scores[$r]
}

# Summary row: average of the Score column.
$ws.Range("D33").Formula = "=AVERAGE(D2:D32)"

# Update sheet view state to match the recorded selection/scroll position.
$ws.Range("A2").Select()
$ws.Range("D34").Select()

# Update workbook-level window position recorded in the file.
$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
